$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 432
$ws.Range("I4").Value = 90
$ws.Range("K4").Value = 90
$ws.Range("M4").Value = 24
$ws.Range("H5").Value = 107.75
$ws.Range("J5").Value = 274.5
$ws.Range("L5").Value = 274.5
$ws.Range("N5").Value = -504.5
$ws.Range("H33").Value = 55736.332
$ws.Range("I33").Value = 83506.5
$ws.Range("K33").Value = 83506.5
$ws.Range("M33").Value = -83277.5
$ws.Range("H40").Value = 6456
$ws.Range("I40").Value = 4999
$ws.Range("J40").Value = 6698.8335
$ws.Range("K40").Value = 4999
$ws.Range("L40").Value = 6698.8335
$ws.Range("M40").Value = -4824
$ws.Range("N40").Value = -7048.8335
$ws.Range("H53").Value = 780.913
$ws.Range("I53").Value = 842.86664
$ws.Range("K53").Value = 842.86664
$ws.Range("M53").Value = -205.86664
$ws.Range("H58").Value = 1754.3334
$ws.Range("I58").Value = 1754.3334
$ws.Range("K58").Value = 5263.0002
$ws.Range("M58").Value = -5113.0002
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = $null
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = $null
$ws.Range("H70").Value = 3480.3333
$ws.Range("I70").Value = 2220.5
$ws.Range("J70").Value = 6000
$ws.Range("K70").Value = 6661.5
$ws.Range("L70").Value = 18000
$ws.Range("M70").Value = -6391.5
$ws.Range("N70").Value = -18540
$ws.Range("H73").Value = 3480.3333
$ws.Range("I73").Value = 2220.5
$ws.Range("J73").Value = 6000
$ws.Range("K73").Value = 6661.5
$ws.Range("L73").Value = 18000
$ws.Range("M73").Value = -5725.5
$ws.Range("N73").Value = -19872
$ws.Range("H96").Value = 1053
$ws.Range("I96").Value = 1053
$ws.Range("K96").Value = 3159
$ws.Range("M96").Value = -1786
$ws.Range("H113").Value = 6159.75
$ws.Range("I113").Value = 5725.5713
$ws.Range("K113").Value = 5725.5713
$ws.Range("M113").Value = -2471.5713
$ws.Range("H132").Value = 2647.5557
$ws.Range("I132").Value = 2709.3076
$ws.Range("J132").Value = 2487
$ws.Range("K132").Value = 8127.9228
$ws.Range("L132").Value = 7461
$ws.Range("M132").Value = -5597.9228
$ws.Range("N132").Value = -12521
$ws.Range("H138").Value = 4299.6875
$ws.Range("J138").Value = 5127.4
$ws.Range("L138").Value = 15382.2
$ws.Range("N138").Value = -25662.2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 452.66666
$ws.Range("I4").Value = 452.66666
$ws.Range("K4").Value = 452.66666
$ws.Range("M4").Value = -336.66666
$ws.Range("H16").Value = 8601
$ws.Range("I16").Value = 8001.25
$ws.Range("J16").Value = 11000
$ws.Range("K16").Value = 8001.25
$ws.Range("L16").Value = 11000
$ws.Range("M16").Value = -7714.25
$ws.Range("N16").Value = -11574
$ws.Range("H25").Value = 4666.6665
$ws.Range("I25").Value = 1700
$ws.Range("K25").Value = 1700
$ws.Range("M25").Value = -1298
$ws.Range("H32").Value = 3094.2432
$ws.Range("I32").Value = 2924.1516
$ws.Range("K32").Value = 2924.1516
$ws.Range("M32").Value = -2637.1516
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").Value = $null
$ws.Range("H45").Value = 797.8
$ws.Range("I45").Value = 612.9231
$ws.Range("K45").Value = 612.9231
$ws.Range("M45").Value = -235.9231
$ws.Range("H74").Value = 2968.2
$ws.Range("I74").Value = 3043.2415
$ws.Range("K74").Value = 3043.2415
$ws.Range("M74").Value = -2169.2415
$ws.Range("H77").Value = 2968.2
$ws.Range("I77").Value = 3043.2415
$ws.Range("K77").Value = 15216.2075
$ws.Range("M77").Value = -10848.2075
$ws.Range("H80").Value = 30110
$ws.Range("J80").Value = 30110
$ws.Range("L80").Value = 30110
$ws.Range("N80").Value = -32106
$ws.Range("H83").Value = 30110
$ws.Range("J83").Value = 30110
$ws.Range("L83").Value = 90330
$ws.Range("N83").Value = -100314
$ws.Range("H88").Value = 2919.6
$ws.Range("I88").Value = 2700
$ws.Range("J88").Value = 2974.5
$ws.Range("K88").Value = 2700
$ws.Range("L88").Value = 2974.5
$ws.Range("M88").Value = -2294
$ws.Range("N88").Value = -3786.5
$ws.Range("H91").Value = 2919.6
$ws.Range("I91").Value = 2700
$ws.Range("J91").Value = 2974.5
$ws.Range("K91").Value = 2700
$ws.Range("L91").Value = 2974.5
$ws.Range("M91").Value = -1296
$ws.Range("N91").Value = -5782.5
$ws.Range("H102").Value = 2481
$ws.Range("I102").Value = 2481
$ws.Range("K102").Value = 2481
$ws.Range("M102").Value = -859
$ws.Range("H109").Value = 70376.5
$ws.Range("J109").Value = 70376.5
$ws.Range("L109").Value = 70376.5
$ws.Range("N109").Value = -73150.5
$ws.Range("H132").Value = 5363.837
$ws.Range("I132").Value = 5329.222
$ws.Range("K132").Value = 15987.666
$ws.Range("M132").Value = -13457.666
$ws.Range("H135").Value = 22999.666
$ws.Range("J135").Value = 22999.666
$ws.Range("L135").Value = 22999.666
$ws.Range("N135").Value = -33139.666
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1724
$ws.Range("I64").Value = 799
$ws.Range("J64").Value = 2649
$ws.Range("K64").Value = 799
$ws.Range("L64").Value = 2649
$ws.Range("M64").Value = -574
$ws.Range("N64").Value = -3099
$ws.Range("H67").Value = 1724
$ws.Range("I67").Value = 799
$ws.Range("J67").Value = 2649
$ws.Range("K67").Value = 799
$ws.Range("L67").Value = 2649
$ws.Range("M67").Value = -19
$ws.Range("N67").Value = -4209
$ws.Range("H86").Value = 180749.75
$ws.Range("I86").Value = 7000
$ws.Range("K86").Value = 7000
$ws.Range("M86").Value = -5877
$ws.Range("H89").Value = 180749.75
$ws.Range("I89").Value = 7000
$ws.Range("K89").Value = 35000
$ws.Range("M89").Value = -29384
$ws.Range("H94").Value = 1247.9231
$ws.Range("I94").Value = 515.3333
$ws.Range("K94").Value = 515.3333
$ws.Range("M94").Value = -64.33330000000001
$ws.Range("H99").Value = 2788
$ws.Range("I99").Value = 2470.6667
$ws.Range("K99").Value = 2470.6667
$ws.Range("M99").Value = -972.6667000000002
$ws.Range("H100").Value = 18150
$ws.Range("J100").Value = 18150
$ws.Range("L100").Value = 18150
$ws.Range("N100").Value = -20314
$ws.Range("H105").Value = 8252.5
$ws.Range("I105").Value = 3999.5
$ws.Range("J105").Value = 12505.5
$ws.Range("K105").Value = 3999.5
$ws.Range("L105").Value = 12505.5
$ws.Range("M105").Value = -2252.5
$ws.Range("N105").Value = -15999.5
$ws.Range("H132").Value = 89889.5
$ws.Range("J132").Value = 89889.5
$ws.Range("L132").Value = 89889.5
$ws.Range("N132").Value = -100009.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 933.26666
$ws.Range("I7").Value = 261.42856
$ws.Range("K7").Value = 261.42856
$ws.Range("M7").Value = -148.42856
$ws.Range("H15").Value = 4332.909
$ws.Range("J15").Value = 4945.8887
$ws.Range("L15").Value = 4945.8887
$ws.Range("N15").Value = -5285.8887
$ws.Range("H25").Value = 14799
$ws.Range("I25").Value = 10598
$ws.Range("J25").Value = 19000
$ws.Range("K25").Value = 10598
$ws.Range("L25").Value = 19000
$ws.Range("M25").Value = -10424
$ws.Range("N25").Value = -19348
$ws.Range("H31").Value = 1501.6666
$ws.Range("I31").Value = 1819.5385
$ws.Range("J31").Value = 1322
$ws.Range("K31").Value = 1819.5385
$ws.Range("L31").Value = 1322
$ws.Range("M31").Value = -1524.5385
$ws.Range("N31").Value = -1912
$ws.Range("H34").Value = 1501.6666
$ws.Range("I34").Value = 1819.5385
$ws.Range("J34").Value = 1322
$ws.Range("K34").Value = 1819.5385
$ws.Range("L34").Value = 1322
$ws.Range("M34").Value = -1617.5385
$ws.Range("N34").Value = -1726
$ws.Range("H58").Value = 2435.6316
$ws.Range("I58").Value = 2216.7334
$ws.Range("K58").Value = 2216.7334
$ws.Range("M58").Value = -2013.7334
$ws.Range("H62").Value = 20647.875
$ws.Range("I62").Value = 6119.25
$ws.Range("J62").Value = 25490.75
$ws.Range("K62").Value = 6119.25
$ws.Range("L62").Value = 25490.75
$ws.Range("M62").Value = -5495.25
$ws.Range("N62").Value = -26738.75
$ws.Range("H65").Value = 20647.875
$ws.Range("I65").Value = 6119.25
$ws.Range("J65").Value = 25490.75
$ws.Range("K65").Value = 30596.25
$ws.Range("L65").Value = 127453.75
$ws.Range("M65").Value = -27476.25
$ws.Range("N65").Value = -133693.75
$ws.Range("H93").Value = 10058.8
$ws.Range("I93").Value = 7249.75
$ws.Range("J93").Value = 21295
$ws.Range("K93").Value = 7249.75
$ws.Range("L93").Value = 21295
$ws.Range("M93").Value = -5377.75
$ws.Range("N93").Value = -25039
$ws.Range("H107").Value = 4253.615
$ws.Range("J107").Value = 7605.75
$ws.Range("L107").Value = 7605.75
$ws.Range("N107").Value = -11445.75
$ws.Range("H135").Value = 199999
$ws.Range("J135").Value = 199999
$ws.Range("L135").Value = 199999
$ws.Range("N135").Value = -210139
$ws.Range("H136").Value = 2435.6316
$ws.Range("I136").Value = 2216.7334
$ws.Range("K136").Value = 6650.2002
$ws.Range("M136").Value = -4100.2002
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9176.666999999999
$ws.Range("I3").Value = 9176.666999999999
$ws.Range("K3").Value = 27530.001
$ws.Range("M3").Value = -27418.001
$ws.Range("H4").Value = 424961.3
$ws.Range("I4").Value = 616077.25
$ws.Range("J4").Value = 233845.39
$ws.Range("K4").Value = 1848231.75
$ws.Range("L4").Value = 701536.17
$ws.Range("M4").Value = -1848119.75
$ws.Range("N4").Value = -701760.17
$ws.Range("H11").Value = 1041.6
$ws.Range("I11").Value = 453.33334
$ws.Range("J11").Value = 1924
$ws.Range("K11").Value = 1360.00002
$ws.Range("L11").Value = 5772
$ws.Range("M11").Value = -1220.00002
$ws.Range("N11").Value = -6052
$ws.Range("H23").Value = 1000.0909
$ws.Range("I23").Value = 803.2727
$ws.Range("J23").Value = 1196.909
$ws.Range("K23").Value = 2409.8181
$ws.Range("L23").Value = 3590.727
$ws.Range("M23").Value = -2174.8181
$ws.Range("N23").Value = -4060.727
$ws.Range("H25").Value = 833.3333
$ws.Range("I25").Value = 750
$ws.Range("J25").Value = 1000
$ws.Range("K25").Value = 2250
$ws.Range("L25").Value = 3000
$ws.Range("M25").Value = -2081
$ws.Range("N25").Value = -3338
$ws.Range("H30").Value = 833.3333
$ws.Range("I30").Value = 750
$ws.Range("J30").Value = 1000
$ws.Range("K30").Value = 2250
$ws.Range("L30").Value = 3000
$ws.Range("M30").Value = -2148
$ws.Range("N30").Value = -3204
$ws.Range("H34").Value = 4609.625
$ws.Range("J34").Value = 8002.3076
$ws.Range("L34").Value = 24006.9228
$ws.Range("N34").Value = -24174.9228
$ws.Range("H39").Value = 7437.091
$ws.Range("J39").Value = 7812
$ws.Range("L39").Value = 23436
$ws.Range("N39").Value = -24024
$ws.Range("H48").Value = 10833
$ws.Range("J48").Value = 10833
$ws.Range("L48").Value = 32499
$ws.Range("N48").Value = -32999
$ws.Range("H55").Value = 4712.1113
$ws.Range("J55").Value = 5772.857
$ws.Range("L55").Value = 17318.571
$ws.Range("N55").Value = -17672.571
$ws.Range("H57").Value = 4500
$ws.Range("J57").Value = 4500
$ws.Range("L57").Value = 13500
$ws.Range("N57").Value = -14618
$ws.Range("H68").Value = 1180.6545
$ws.Range("I68").Value = 905.85297
$ws.Range("J68").Value = 1625.5714
$ws.Range("K68").Value = 2717.55891
$ws.Range("L68").Value = 4876.7142
$ws.Range("M68").Value = -1906.55891
$ws.Range("N68").Value = -6498.7142
$ws.Range("H71").Value = 1180.6545
$ws.Range("I71").Value = 905.85297
$ws.Range("J71").Value = 1625.5714
$ws.Range("K71").Value = 8152.67673
$ws.Range("L71").Value = 14630.1426
$ws.Range("M71").Value = -4096.67673
$ws.Range("N71").Value = -22742.1426
$ws.Range("H107").Value = 1838.3636
$ws.Range("J107").Value = 2002.25
$ws.Range("L107").Value = 6006.75
$ws.Range("N107").Value = -9846.75
$ws.Range("H121").Value = 51131.5
$ws.Range("I121").Value = 148032.42
$ws.Range("J121").Value = 5911.067
$ws.Range("K121").Value = 444097.26
$ws.Range("L121").Value = 17733.201
$ws.Range("M121").Value = -442787.26
$ws.Range("N121").Value = -20353.201
$ws.Range("H138").Value = 8166.6665
$ws.Range("I138").Value = 7250
$ws.Range("J138").Value = 10000
$ws.Range("K138").Value = 21750
$ws.Range("L138").Value = 30000
$ws.Range("M138").Value = -16610
$ws.Range("N138").Value = -40280
$ws.Range("H139").Value = 1646.6666
$ws.Range("I139").Value = 1352.75
$ws.Range("K139").Value = 4058.25
$ws.Range("M139").Value = 1081.75
$ws.Range("H140").Value = 2196.1875
$ws.Range("I140").Value = 1219.8846
$ws.Range("J140").Value = 6426.8335
$ws.Range("K140").Value = 3659.6538
$ws.Range("L140").Value = 19280.5005
$ws.Range("M140").Value = 1520.3462
$ws.Range("N140").Value = -29640.5005
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 10336888
$ws.Range("I14").Value = 15501332
$ws.Range("K14").Value = 15501332
$ws.Range("M14").Value = -15501164
$ws.Range("H19").Value = 9265.143
$ws.Range("J19").Value = 9265.143
$ws.Range("L19").Value = 9265.143
$ws.Range("N19").Value = -9841.143
$ws.Range("H28").Value = 126015
$ws.Range("J28").Value = 126015
$ws.Range("L28").Value = 126015
$ws.Range("N28").Value = -126399
$ws.Range("H69").Value = 55191.5
$ws.Range("J69").Value = 60201
$ws.Range("L69").Value = 60201
$ws.Range("N69").Value = -61699
$ws.Range("H72").Value = 55191.5
$ws.Range("J72").Value = 60201
$ws.Range("L72").Value = 180603
$ws.Range("N72").Value = -188091
$ws.Range("H122").Value = 2397
$ws.Range("I122").Value = 2397
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7191
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4741
$ws.Range("N122").Value = $null
$ws.Range("H126").Value = 3272.8572
$ws.Range("I126").Value = 2970
$ws.Range("K126").Value = 8910
$ws.Range("M126").Value = -6440
$ws.Range("H132").Value = 2950.0688
$ws.Range("I132").Value = 2831.9565
$ws.Range("J132").Value = 3402.8333
$ws.Range("K132").Value = 8495.869499999999
$ws.Range("L132").Value = 10208.4999
$ws.Range("M132").Value = -5965.869499999999
$ws.Range("N132").Value = -15268.4999
$ws.Range("H135").Value = 695475
$ws.Range("J135").Value = 695475
$ws.Range("L135").Value = 695475
$ws.Range("N135").Value = -705615
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1051.5
$ws.Range("I7").Value = 1051.5
$ws.Range("K7").Value = 1051.5
$ws.Range("M7").Value = -939.5
$ws.Range("H19").Value = 21334.666
$ws.Range("J19").Value = 60004
$ws.Range("L19").Value = 60004
$ws.Range("N19").Value = -60344
$ws.Range("H40").Value = 3971.2
$ws.Range("I40").Value = 3266.8462
$ws.Range("J40").Value = 8549.5
$ws.Range("K40").Value = 3266.8462
$ws.Range("L40").Value = 8549.5
$ws.Range("M40").Value = -3130.8462
$ws.Range("N40").Value = -8821.5
$ws.Range("H61").Value = 1685.1428
$ws.Range("I61").Value = 1685.1428
$ws.Range("K61").Value = 1685.1428
$ws.Range("M61").Value = -1483.1428
$ws.Range("H93").Value = 3718.4
$ws.Range("I93").Value = 3491.8
$ws.Range("J93").Value = 4058.3
$ws.Range("K93").Value = 3491.8
$ws.Range("L93").Value = 4058.3
$ws.Range("M93").Value = -2243.8
$ws.Range("N93").Value = -6554.3
$ws.Range("H113").Value = 1685.1428
$ws.Range("I113").Value = 1685.1428
$ws.Range("K113").Value = 1685.1428
$ws.Range("M113").Value = 484.8571999999999
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = $null
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = $null
$ws.Range("H126").Value = 1051.5
$ws.Range("I126").Value = 1051.5
$ws.Range("K126").Value = 3154.5
$ws.Range("M126").Value = -684.5
$ws.Range("H132").Value = 2754
$ws.Range("I132").Value = 2503
$ws.Range("K132").Value = 7509
$ws.Range("M132").Value = -4979
$ws.Range("H136").Value = 90915544
$ws.Range("I136").Value = 6933.3335
$ws.Range("K136").Value = 20800.0005
$ws.Range("M136").Value = -18250.0005
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 27601
$ws.Range("J14").Value = 27601
$ws.Range("L14").Value = 27601
$ws.Range("N14").Value = -27937
$ws.Range("H17").Value = 24500
$ws.Range("I17").Value = 17000
$ws.Range("J17").Value = 32000
$ws.Range("K17").Value = 17000
$ws.Range("L17").Value = 32000
$ws.Range("M17").Value = -16828
$ws.Range("N17").Value = -32344
$ws.Range("H74").Value = 12313
$ws.Range("J74").Value = 12852.25
$ws.Range("L74").Value = 12852.25
$ws.Range("N74").Value = -14724.25
$ws.Range("H77").Value = 12313
$ws.Range("J77").Value = 12852.25
$ws.Range("L77").Value = 38556.75
$ws.Range("N77").Value = -47916.75
$ws.Range("H88").Value = 20000
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").Value = $null
$ws.Range("H91").Value = 20000
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").Value = $null
$ws.Range("H122").Value = 2292.7
$ws.Range("I122").Value = 2292.7
$ws.Range("K122").Value = 6878.099999999999
$ws.Range("M122").Value = -4428.099999999999
$ws.Range("H126").Value = 2760.35
$ws.Range("I126").Value = 2233.7778
$ws.Range("J126").Value = 7499.5
$ws.Range("K126").Value = 6701.3334
$ws.Range("L126").Value = 22498.5
$ws.Range("M126").Value = -4231.3334
$ws.Range("N126").Value = -27438.5
$ws.Range("H132").Value = 2857.5557
$ws.Range("I132").Value = 2589.75
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 7769.25
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -5239.25
$ws.Range("N132").Value = -20060
$ws.Range("H135").Value = 60000
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140
$ws.Range("H136").Value = 1879.2593
$ws.Range("I136").Value = 1669.64
$ws.Range("K136").Value = 5008.92
$ws.Range("M136").Value = -2458.92
